$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44335
$ws.Cells.Item(2, 10).Value = 250
$ws.Cells.Item(2, 11).Value = 1400
$ws.Cells.Item(2, 12).Value = 1500
$ws.Cells.Item(2, 13).Value = 1450
$ws.Cells.Item(2, 16).Value = 483
$ws.Cells.Item(3, 4).Value = 44432
$ws.Cells.Item(3, 10).Value = 270
$ws.Cells.Item(3, 11).Value = 1800
$ws.Cells.Item(3, 12).Value = 2000
$ws.Cells.Item(3, 13).Value = 1900
$ws.Cells.Item(3, 16).Value = 633
$ws.Cells.Item(4, 4).Value = 44412
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 2800
$ws.Cells.Item(4, 12).Value = 3000
$ws.Cells.Item(4, 13).Value = 2900
$ws.Cells.Item(4, 16).Value = 967
$ws.Cells.Item(5, 4).Value = 44364
$ws.Cells.Item(5, 10).Value = 270
$ws.Cells.Item(5, 11).Value = 3400
$ws.Cells.Item(5, 12).Value = 3500
$ws.Cells.Item(5, 13).Value = 3450
$ws.Cells.Item(5, 16).Value = 1150
$ws.Cells.Item(6, 4).Value = 44435
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 270
$ws.Cells.Item(6, 11).Value = 1800
$ws.Cells.Item(6, 12).Value = 2000
$ws.Cells.Item(6, 13).Value = 1900
$ws.Cells.Item(6, 16).Value = 633
$ws.Cells.Item(7, 4).Value = 44313
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 11).Value = 1300
$ws.Cells.Item(7, 12).Value = 1500
$ws.Cells.Item(7, 13).Value = 1400
$ws.Cells.Item(7, 16).Value = 467
$ws.Cells.Item(8, 4).Value = 44313
$ws.Cells.Item(8, 9).Value = 'Segunda'
$ws.Cells.Item(8, 10).Value = 250
$ws.Cells.Item(8, 11).Value = 900
$ws.Cells.Item(8, 12).Value = 1000
$ws.Cells.Item(8, 13).Value = 950
$ws.Cells.Item(8, 16).Value = 317
$ws.Cells.Item(9, 4).Value = 44322
$ws.Cells.Item(9, 10).Value = 250
$ws.Cells.Item(10, 4).Value = 44406
$ws.Cells.Item(10, 11).Value = 2800
$ws.Cells.Item(10, 12).Value = 3000
$ws.Cells.Item(10, 13).Value = 2900
$ws.Cells.Item(10, 16).Value = 967
$ws.Cells.Item(11, 4).Value = 44299
$ws.Cells.Item(11, 10).Value = 300
$ws.Cells.Item(11, 11).Value = 1400
$ws.Cells.Item(11, 13).Value = 1450
$ws.Cells.Item(11, 16).Value = 483
$ws.Cells.Item(12, 4).Value = 44299
$ws.Cells.Item(12, 10).Value = 250
$ws.Cells.Item(12, 11).Value = 1000
$ws.Cells.Item(12, 12).Value = 1200
$ws.Cells.Item(12, 13).Value = 1100
$ws.Cells.Item(12, 16).Value = 367
$ws.Cells.Item(15, 4).Value = 44327
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 1400
$ws.Cells.Item(15, 12).Value = 1500
$ws.Cells.Item(15, 13).Value = 1450
$ws.Cells.Item(15, 16).Value = 483
$ws.Cells.Item(16, 4).Value = 44327
$ws.Cells.Item(16, 9).Value = 'Segunda'
$ws.Cells.Item(16, 10).Value = 250
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 12).Value = 1200
$ws.Cells.Item(16, 13).Value = 1100
$ws.Cells.Item(16, 16).Value = 367
$ws.Cells.Item(17, 4).Value = 44300
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 160
$ws.Cells.Item(17, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(18, 4).Value = 44383
$ws.Cells.Item(18, 9).Value = 'Segunda'
$ws.Cells.Item(18, 10).Value = 350
$ws.Cells.Item(18, 11).Value = 2800
$ws.Cells.Item(18, 12).Value = 3000
$ws.Cells.Item(18, 13).Value = 2886
$ws.Cells.Item(18, 16).Value = 962
$ws.Cells.Item(19, 4).Value = 44341
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 1400
$ws.Cells.Item(19, 12).Value = 1500
$ws.Cells.Item(19, 13).Value = 1450
$ws.Cells.Item(19, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(19, 16).Value = 483
$ws.Cells.Item(20, 4).Value = 44392
$ws.Cells.Item(20, 10).Value = 200
$ws.Cells.Item(20, 11).Value = 3800
$ws.Cells.Item(20, 12).Value = 4000
$ws.Cells.Item(20, 13).Value = 3900
$ws.Cells.Item(20, 16).Value = 1300
$ws.Cells.Item(21, 9).Value = 'Segunda'
$ws.Cells.Item(21, 11).Value = 3200
$ws.Cells.Item(21, 12).Value = 3500
$ws.Cells.Item(21, 13).Value = 3350
$ws.Cells.Item(21, 16).Value = 1117
$ws.Cells.Item(22, 4).Value = 44448
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 11).Value = 1400
$ws.Cells.Item(22, 12).Value = 1500
$ws.Cells.Item(22, 13).Value = 1450
$ws.Cells.Item(22, 16).Value = 483
$ws.Cells.Item(23, 4).Value = 44448
$ws.Cells.Item(23, 9).Value = 'Segunda'
$ws.Cells.Item(23, 10).Value = 200
$ws.Cells.Item(23, 11).Value = 1000
$ws.Cells.Item(23, 12).Value = 1200
$ws.Cells.Item(23, 13).Value = 1100
$ws.Cells.Item(23, 16).Value = 367
$ws.Cells.Item(24, 4).Value = 44217
$ws.Cells.Item(24, 10).Value = 250
$ws.Cells.Item(24, 11).Value = 2900
$ws.Cells.Item(24, 12).Value = 3000
$ws.Cells.Item(24, 13).Value = 2950
$ws.Cells.Item(24, 16).Value = 983
$ws.Cells.Item(25, 4).Value = 44350
$ws.Cells.Item(25, 10).Value = 300
$ws.Cells.Item(25, 11).Value = 1800
$ws.Cells.Item(25, 12).Value = 2000
$ws.Cells.Item(25, 13).Value = 1900
$ws.Cells.Item(25, 16).Value = 633
$ws.Cells.Item(26, 4).Value = 44418
$ws.Cells.Item(26, 10).Value = 300
$ws.Cells.Item(26, 11).Value = 2400
$ws.Cells.Item(26, 12).Value = 2500
$ws.Cells.Item(26, 13).Value = 2450
$ws.Cells.Item(26, 16).Value = 817
$ws.Cells.Item(27, 4).Value = 44343
$ws.Cells.Item(27, 10).Value = 150
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 13).Value = 1500
$ws.Cells.Item(27, 16).Value = 500
$ws.Cells.Item(28, 4).Value = 44343
$ws.Cells.Item(28, 9).Value = 'Segunda'
$ws.Cells.Item(28, 10).Value = 150
$ws.Cells.Item(28, 11).Value = 1400
$ws.Cells.Item(28, 12).Value = 1400
$ws.Cells.Item(28, 13).Value = 1400
$ws.Cells.Item(28, 16).Value = 467
$ws.Cells.Item(29, 4).Value = 44308
$ws.Cells.Item(29, 10).Value = 270
$ws.Cells.Item(29, 11).Value = 1400
$ws.Cells.Item(29, 12).Value = 1500
$ws.Cells.Item(29, 13).Value = 1450
$ws.Cells.Item(29, 16).Value = 483
$ws.Cells.Item(30, 4).Value = 44168
$ws.Cells.Item(30, 10).Value = 300
$ws.Cells.Item(30, 11).Value = 1800
$ws.Cells.Item(30, 12).Value = 2000
$ws.Cells.Item(30, 13).Value = 1900
$ws.Cells.Item(30, 16).Value = 633
$ws.Cells.Item(31, 4).Value = 44277
$ws.Cells.Item(31, 10).Value = 250
$ws.Cells.Item(31, 11).Value = 1800
$ws.Cells.Item(31, 12).Value = 2000
$ws.Cells.Item(31, 13).Value = 1900
$ws.Cells.Item(31, 16).Value = 633
$ws.Cells.Item(32, 4).Value = 44273
$ws.Cells.Item(32, 9).Value = 'Primera'
$ws.Cells.Item(32, 11).Value = 3800
$ws.Cells.Item(32, 12).Value = 4000
$ws.Cells.Item(32, 13).Value = 3900
$ws.Cells.Item(32, 16).Value = 1300
$ws.Cells.Item(33, 4).Value = 44356
$ws.Cells.Item(33, 10).Value = 200
$ws.Cells.Item(33, 11).Value = 2400
$ws.Cells.Item(33, 12).Value = 2500
$ws.Cells.Item(33, 13).Value = 2450
$ws.Cells.Item(33, 16).Value = 817
$ws.Cells.Item(34, 4).Value = 44356
$ws.Cells.Item(34, 9).Value = 'Segunda'
$ws.Cells.Item(34, 10).Value = 200
$ws.Cells.Item(35, 4).Value = 44175
$ws.Cells.Item(35, 10).Value = 250
$ws.Cells.Item(35, 11).Value = 1800
$ws.Cells.Item(35, 12).Value = 2000
$ws.Cells.Item(35, 13).Value = 1900
$ws.Cells.Item(35, 16).Value = 633
$ws.Cells.Item(36, 4).Value = 44257
$ws.Cells.Item(36, 10).Value = 1500
$ws.Cells.Item(37, 4).Value = 44376
$ws.Cells.Item(37, 10).Value = 280
$ws.Cells.Item(37, 11).Value = 2400
$ws.Cells.Item(37, 12).Value = 2500
$ws.Cells.Item(37, 13).Value = 2436
$ws.Cells.Item(37, 16).Value = 812
$ws.Cells.Item(38, 4).Value = 44292
$ws.Cells.Item(38, 9).Value = 'Primera'
$ws.Cells.Item(38, 10).Value = 270
$ws.Cells.Item(38, 11).Value = 2400
$ws.Cells.Item(38, 12).Value = 2500
$ws.Cells.Item(38, 13).Value = 2450
$ws.Cells.Item(38, 16).Value = 817
